# Scheduled data-refresh runner: re-pull Universalis market-board prices
# for each job sheet and rewrite the derived Leve price/profit columns
# (H:N) with the freshly computed figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1347.75
$ws.Range("I41").Value = 1725.6
$ws.Range("J41").Value = 718
$ws.Range("K41").Value = 1725.6
$ws.Range("L41").Value = 718
$ws.Range("M41").Value = -1285.6
$ws.Range("N41").Value = -1598

$ws.Range("H69").Value = 12767.929
$ws.Range("J69").Value = 12106.571
$ws.Range("L69").Value = 36319.713
$ws.Range("N69").Value = -38067.713

$ws.Range("H72").Value = 12767.929
$ws.Range("J72").Value = 12106.571
$ws.Range("L72").Value = 108959.139
$ws.Range("N72").Value = -117695.139

$ws.Range("H125").Value = 4212.7334
$ws.Range("I125").Value = 3734.8
$ws.Range("J125").Value = 4451.7
$ws.Range("K125").Value = 33613.2
$ws.Range("L125").Value = 40065.3
$ws.Range("M125").Value = -31153.2
$ws.Range("N125").Value = -44985.3

$ws.Range("H132").Value = 2071.6667
$ws.Range("I132").Value = 2071.4365
$ws.Range("J132").Value = 2074
$ws.Range("K132").Value = 6214.309499999999
$ws.Range("L132").Value = 6222
$ws.Range("M132").Value = -3684.309499999999
$ws.Range("N132").Value = -11282

$ws.Range("H138").Value = 2980.5789
$ws.Range("I138").Value = 2799.5454
$ws.Range("J138").Value = 3229.5
$ws.Range("K138").Value = 8398.636200000001
$ws.Range("L138").Value = 9688.5
$ws.Range("M138").Value = -3258.636200000001
$ws.Range("N138").Value = -19968.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1259.2307
$ws.Range("I2").Value = 947.5
$ws.Range("K2").Value = 947.5
$ws.Range("M2").Value = -834.5

$ws.Range("H32").Value = 4156.149
$ws.Range("I32").Value = 4551.9736
$ws.Range("K32").Value = 4551.9736
$ws.Range("M32").Value = -4264.9736

$ws.Range("H45").Value = 1680.4762
$ws.Range("I45").Value = 1066
$ws.Range("J45").Value = 3646.8
$ws.Range("K45").Value = 1066
$ws.Range("L45").Value = 3646.8
$ws.Range("M45").Value = -689
$ws.Range("N45").Value = -4400.8

$ws.Range("H61").Value = 3569.9119
$ws.Range("I61").Value = 3120.2144
$ws.Range("J61").Value = 5668.5
$ws.Range("K61").Value = 3120.2144
$ws.Range("L61").Value = 5668.5
$ws.Range("M61").Value = -2908.2144
$ws.Range("N61").Value = -6092.5

$ws.Range("H74").Value = 1811
$ws.Range("I74").Value = 1743.8334
$ws.Range("J74").Value = 2214
$ws.Range("K74").Value = 1743.8334
$ws.Range("L74").Value = 2214
$ws.Range("M74").Value = -869.8334
$ws.Range("N74").Value = -3962

$ws.Range("H77").Value = 1811
$ws.Range("I77").Value = 1743.8334
$ws.Range("J77").Value = 2214
$ws.Range("K77").Value = 8719.166999999999
$ws.Range("L77").Value = 11070
$ws.Range("M77").Value = -4351.166999999999
$ws.Range("N77").Value = -19806

$ws.Range("H116").Value = 1259.2307
$ws.Range("I116").Value = 947.5
$ws.Range("K116").Value = 947.5
$ws.Range("M116").Value = 1346.5

$ws.Range("H136").Value = 3569.9119
$ws.Range("I136").Value = 3120.2144
$ws.Range("J136").Value = 5668.5
$ws.Range("K136").Value = 9360.643199999999
$ws.Range("L136").Value = 17005.5
$ws.Range("M136").Value = -6810.643199999999
$ws.Range("N136").Value = -22105.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1259.2307
$ws.Range("I3").Value = 947.5
$ws.Range("K3").Value = 947.5
$ws.Range("M3").Value = -833.5

$ws.Range("H94").Value = 14493777
$ws.Range("I94").Value = 23810298
$ws.Range("K94").Value = 23810298
$ws.Range("M94").Value = -23809847

$ws.Range("H134").Value = 4632.5747
$ws.Range("I134").Value = 4254.575
$ws.Range("K134").Value = 12763.725
$ws.Range("M134").Value = -10228.725

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 38486.168
$ws.Range("J60").Value = 38486.168
$ws.Range("L60").Value = 38486.168
$ws.Range("N60").Value = -39508.168

$ws.Range("H99").Value = 2154.3333
$ws.Range("I99").Value = 2266.913
$ws.Range("K99").Value = 2266.913
$ws.Range("M99").Value = -768.913

$ws.Range("H107").Value = 1688.762
$ws.Range("I107").Value = 832.7857
$ws.Range("K107").Value = 832.7857
$ws.Range("M107").Value = 1087.2143

$ws.Range("H126").Value = 2154.3333
$ws.Range("I126").Value = 2266.913
$ws.Range("K126").Value = 6800.739
$ws.Range("M126").Value = -4330.739

$ws.Range("H127").Value = 69926.664
$ws.Range("J127").Value = 69926.664
$ws.Range("L127").Value = 69926.664
$ws.Range("N127").Value = -79846.664

$ws.Range("H134").Value = 2366.5151
$ws.Range("I134").Value = 2416
$ws.Range("J134").Value = 1599.5
$ws.Range("K134").Value = 7248
$ws.Range("L134").Value = 4798.5
$ws.Range("M134").Value = -4713
$ws.Range("N134").Value = -9868.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12525833
$ws.Range("I4").Value = 16410596
$ws.Range("J4").Value = 53698.21
$ws.Range("K4").Value = 49231788
$ws.Range("L4").Value = 161094.63
$ws.Range("M4").Value = -49231676
$ws.Range("N4").Value = -161318.63

$ws.Range("H5").Value = 680.907
$ws.Range("I5").Value = 369.43478
$ws.Range("K5").Value = 1108.30434
$ws.Range("M5").Value = -996.3043399999999

$ws.Range("H32").Value = 983.75
$ws.Range("J32").Value = 983.75
$ws.Range("L32").Value = 2951.25
$ws.Range("N32").Value = -3517.25

$ws.Range("H68").Value = 1924.4828
$ws.Range("J68").Value = 1856.0212
$ws.Range("L68").Value = 5568.063599999999
$ws.Range("N68").Value = -7190.063599999999

$ws.Range("H71").Value = 1924.4828
$ws.Range("J71").Value = 1856.0212
$ws.Range("L71").Value = 16704.1908
$ws.Range("N71").Value = -24816.1908

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H106").Value = 25029
$ws.Range("J106").Value = 25029
$ws.Range("L106").Value = 75087
$ws.Range("N106").Value = -76979

$ws.Range("H135").Value = 680.907
$ws.Range("I135").Value = 369.43478
$ws.Range("K135").Value = 3324.91302
$ws.Range("M135").Value = -789.91302

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1606.5385
$ws.Range("J46").Value = 1606.8334
$ws.Range("L46").Value = 1606.8334
$ws.Range("N46").Value = -1982.8334

$ws.Range("H132").Value = 10201.5
$ws.Range("I132").Value = 9108.177
$ws.Range("J132").Value = 12856.714
$ws.Range("K132").Value = 27324.531
$ws.Range("L132").Value = 38570.142
$ws.Range("M132").Value = -24794.531
$ws.Range("N132").Value = -43630.142

$ws.Range("H136").Value = 8775345
$ws.Range("I136").Value = 2702.2415
$ws.Range("K136").Value = 8106.7245
$ws.Range("M136").Value = -5556.7245

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 900.7027
$ws.Range("I113").Value = 860.0645
$ws.Range("K113").Value = 2580.1935
$ws.Range("M113").Value = -410.1934999999999

$ws.Range("H126").Value = 53887510
$ws.Range("I126").Value = 63990750
$ws.Range("K126").Value = 191972250
$ws.Range("M126").Value = -191969780

$ws.Range("H128").Value = 69714.836
$ws.Range("J128").Value = 69714.836
$ws.Range("L128").Value = 69714.836
$ws.Range("N128").Value = -79674.836

$ws.Range("H129").Value = 69238.336
$ws.Range("J129").Value = 69238.336
$ws.Range("L129").Value = 69238.336
$ws.Range("N129").Value = -79238.336

$ws.Range("H132").Value = 5123.775
$ws.Range("I132").Value = 3476.611
$ws.Range("J132").Value = 19948.25
$ws.Range("K132").Value = 10429.833
$ws.Range("L132").Value = 59844.75
$ws.Range("M132").Value = -7899.832999999999
$ws.Range("N132").Value = -64904.75

$ws.Range("H136").Value = 2313.1614
$ws.Range("I136").Value = 1922.3158
$ws.Range("K136").Value = 5766.9474
$ws.Range("M136").Value = -3216.9474
